$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill the empty/placeholder predicted price cell (TimeTaken in Hours)
# with a real formula derived from the minutes column instead of a
# hard-coded literal value.
$ws.Range("C2").Formula = "=B2/60"
